# SCD0025-001 test-case workbook update
# (commit: "Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025")
#
#  - Rename the worksheet tab SCD0331 -> SCD0025
#  - Update the TC_ID cell (B2) DGS-346 -> SCD0025-001
#  - Widen column B so the longer TC_ID text keeps fitting (best-fit width)
#  - Move the selection/view back to B3 (scrolled back to show column A)
#  - Let the sheet recalc so the TODAY()-driven helper text in L2 refreshes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "SCD0025"

# Update the TC_ID value
$ws.Range("B2").Value = "SCD0025-001"

# Column B needs to be a bit wider to fit the longer TC_ID text
$ws.Columns("B").ColumnWidth = 11.65

# Move the selection to B3 (also resets the scrolled-away top-left cell)
$ws.Range("B3").Select() | Out-Null

# Recalculate so the TODAY()-driven helper text in L2 reflects the current date
$excel.Calculate() | Out-Null
